# PrecioFrutaHortalizas / Hortaliza, Terminal Hortofrutícola Agro Chillán - Haba
# "Fruta / hortaliza, semanal" — insert a new weekly price record as the new
# first data row of the Chillán/Ñuble "Haba" block (old row 38), pushing the
# existing rows 38-60 down to 39-61 and growing the used range to A1:R61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 38; this shifts rows 38:60 down to 39:61
# (same behaviour as right-clicking the row header -> "Insert" in Excel).
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value = "Ñuble"
$ws.Cells.Item(38, 4).Value = 44873
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = 100112026
$ws.Cells.Item(38, 7).Value = "Haba"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 80
$ws.Cells.Item(38, 11).Value = 6500
$ws.Cells.Item(38, 12).Value = 7000
$ws.Cells.Item(38, 13).Value = 6750
$ws.Cells.Item(38, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(38, 16).Value = 270
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"
